$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price (D) and 1h volume change (E) values to the latest scraped figures.
# D-column price cells are forced to Text format before assignment so that values such as
# "1.00", "12.80", "47.320.00" keep their exact original formatting instead of being
# auto-converted into numbers by Excel (which would drop trailing zeros / thousands dots).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "47.320.00"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.38%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.500.05"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.85%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "323.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.02"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.50%  "
$ws.Range("E7").Value = "  +1.29%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -0.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.16"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.61%  "
$ws.Range("E11").Value = "  +0.62%  "
$ws.Range("E12").Value = "  +0.70%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.34"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.18"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.35%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.892.40"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.98%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.504.22"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.45%  "
$ws.Range("E17").Value = "  +1.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "47.250.65"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.57%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.80"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.75%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.64"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0941"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.70"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +12.88%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.59"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.51%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "247.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("E25").Value = "  +3.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.05"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.31%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("E28").Value = "  +0.57%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.06"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.72%  "
$ws.Range("E30").Value = "  +2.41%  "
$ws.Range("E31").Value = "  +7.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.84"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.93%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.98"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.83%  "
$ws.Range("E34").Value = "  +1.64%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0789"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.22%  "
$ws.Range("E36").Value = "  +0.30%  "
$ws.Range("E37").Value = "  +4.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.69"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.16%  "
$ws.Range("E39").Value = "  +1.13%  "
$ws.Range("E40").Value = "  +1.10%  "
$ws.Range("E41").Value = "  +0.52%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "121.21"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.92%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.29"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.60%  "
$ws.Range("E44").Value = "  +2.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.990.44"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.94%  "
$ws.Range("E46").Value = "  +2.83%  "
$ws.Range("E47").Value = "  -1.34%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.79"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.89%  "
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("E50").Value = "  +3.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "56.54"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.81%  "
